# Update loading_percent values for the 380 kV case (Case_3_252).
# Rows 2-25 correspond to data index 0-23 (column A), and the updated
# columns are B, C, D, F, G, I, J, K (E, H, L, M, N, O are unchanged, left as 0).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("B","C","D","F","G","I","J","K")

$data = @(
        @(12.49760842841439,10.31311485474843,6.952287733493436,49.05179531031379,3.731833816890445,32.79219347728996,11.57120687934617,13.92714429132187),
        @(12.36815436899244,10.17389818579716,6.95789911822943,48.75822040275949,3.735061714448932,32.66522112047451,11.55968218319675,13.84357707147129),
        @(12.29244730342642,10.09081884668385,6.96372530009569,48.5857898570913,3.73714579118933,32.59129457117657,11.55500796382237,13.79628871250913),
        @(12.26258531599481,10.05760990460859,6.966698154444735,48.51753304315556,3.738020847014389,32.56219479329508,11.55370812733184,13.77804815087575),
        @(12.25768760564685,10.05213586392019,6.967227942308185,48.50632162171384,3.73816770906511,32.55742506449046,11.55352884751589,13.77508202117208),
        @(12.29204051867663,10.09036830564446,6.963762969530864,48.58486113012662,3.737157488000071,32.59089795373105,11.55498798342117,13.79603852123921),
        @(12.45221266600198,10.26464098668703,6.953728514025712,48.94897275177786,3.732925655276739,32.74758303462028,11.56673531318657,13.89750707436058),
        @(12.79446508917625,10.62348497829853,6.952924896379267,49.72309350978052,3.725433096568447,33.08631831306225,11.60877329974384,14.12750945212752),
        @(13.06052414946532,10.89488718851029,6.963792355107554,50.32562240902146,3.720413626236144,33.35359075913089,11.65114290635816,14.31409950750622),
        @(13.1841665925265,11.01947094925209,6.971204989577963,50.60640272176905,3.718234222255561,33.47898706161649,11.67288204482321,14.40250049878202),
        @(13.23131413534141,11.06676186724516,6.974364665737866,50.71362922305364,3.717423791204032,33.52700185717157,11.68146533057408,14.43645437798521),
        @(13.22114630574065,11.05657260960007,6.973668530032066,50.6904969827995,3.717597672471924,33.51663774430563,11.67960120623018,14.42912098134369),
        @(13.18803917727782,11.02335959073388,6.971457885051721,50.61520660405507,3.718167250264012,33.48292675526734,11.67358117986071,14.40528448096802),
        @(13.16780122640766,11.00302907648275,6.970149651902227,50.56920459643391,3.718518065906883,33.46234625485648,11.66993936407402,14.3907453880468),
        @(13.05249220771999,10.88676383280966,6.963357412109707,50.30740218461885,3.720558140015553,33.34547099221835,11.64977153021478,14.30839077052257),
        @(12.98238799944313,10.81569204941791,6.959821414728617,50.14846512082453,3.721836226546405,33.274735593183,11.63802833036462,14.25875097658032),
        @(12.94231458083364,10.77492180398695,6.958020163526932,50.05768187334601,3.722581140396724,33.23441049280829,11.63150619379609,14.23053255503894),
        @(12.92879059887796,10.76113777291614,6.95745029815573,50.02705478173739,3.722835040008612,33.22081948526925,11.62933789247842,14.22103632175925),
        @(12.98982532804436,10.8232468938282,6.960173774273162,50.16531914682303,3.721699159246499,33.28222836541465,11.63925440069329,14.26400094437159),
        @(13.1977550604198,11.03311234647447,6.972097654648338,50.63729724873192,3.717999548826484,33.49281425818544,11.67533990597899,14.41227308338157),
        @(13.33553175188407,11.17091325082801,6.981944839831887,50.95098523562519,3.715668224818828,33.63352405737741,11.70096879399893,14.51195212220973),
        @(13.26184151120706,11.09732304823376,6.976502142289633,50.78310656796942,3.716904603387289,33.55814908149711,11.68710426333221,14.45850709385915),
        @(12.98646219033377,10.81983106511595,6.960013750818886,50.157697593658,3.721761095827023,33.27883981752483,11.63869938018006,14.261626434568),
        @(12.69912170393811,10.5248515057762,6.951124982026238,49.50752335731908,3.727374368740752,32.99138835739185,11.5953752994632,14.06209544266976)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($j = 0; $j -lt $columns.Count; $j++) {
        $ws.Range("$($columns[$j])$row").Value = $values[$j]
    }
}
